# Horarios actualizados Linea 141 - 432
# Updates the scraped-schedule workbook with the newer scrape snapshot
# (Ultima actualizacion: 08:47:19 -> 08:54:42), which changes the
# "Total filas" counters and appends / reorders several detail rows on
# each of the three sheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 08:54:42"
$ws1.Range("A3").Value = "Total filas: 132"

$sheet1Rows = @(
    @(54,  "06:37:24", "07:16", "16_SANTA ANA",            39,  "LP1912"),
    @(55,  "06:02:16", "07:16", "11_ETCHEVERRY",           74,  "LP1912"),
    @(111, "07:57:27", "09:23", "11_ETCHEVERRY",           86,  "LP1912"),
    @(112, "08:47:19", "09:23", "16_SANTA ANA",             36, "LP1912"),
    @(113, "07:44:08", "09:23", "17_ROMERO",                99, "LP1912"),
    @(114, "07:44:08", "09:24", "11_ETCHEVERRY",           100, "LP1912"),
    @(115, "08:16:48", "09:29", "16_SANTA ANA",             73, "LP1912"),
    @(116, "07:44:08", "09:32", "15_ABASTO",               108, "LP1912"),
    @(117, "07:44:08", "09:33", "10_OLMOS",                109, "LP1912"),
    @(118, "08:33:47", "09:34", "16_SANTA ANA",             61, "LP1912"),
    @(119, "08:54:42", "09:34", "23_HERNANDEZ",             40, "LP1912"),
    @(120, "08:47:19", "09:35", "16_SANTA ANA",             48, "LP1912"),
    @(121, "08:47:19", "09:35", "23_HERNANDEZ",             48, "LP1912"),
    @(122, "07:44:08", "09:36", "23_HERNANDEZ",            112, "LP1912"),
    @(123, "08:16:48", "09:37", "23_HERNANDEZ",             81, "LP1912"),
    @(124, "08:16:48", "09:41", "215C_EL PATO",             85, "LP1912"),
    @(125, "08:33:47", "09:41", "23_HERNANDEZ",             68, "LP1912"),
    @(126, "07:44:08", "09:42", "215C_EL PATO",            118, "LP1912"),
    @(127, "07:57:27", "09:43", "14_ABASTO",                106, "LP1912"),
    @(128, "08:54:42", "09:52", "15_ABASTO",                58, "LP1912"),
    @(129, "08:54:42", "09:53", "10_OLMOS",                 59, "LP1912"),
    @(130, "08:16:48", "10:10", "16_P MOR-SANTA ANA",      114, "LP1912"),
    @(131, "08:16:48", "10:12", "15_ABASTO",                116, "LP1912"),
    @(132, "08:33:47", "10:21", "26_HERNANDEZ",            108, "LP1912"),
    @(133, "08:33:47", "10:22", "17_ROMERO",                109, "LP1912"),
    @(134, "08:33:47", "10:26", "215A_EL PATO",             113, "LP1912"),
    @(135, "08:54:42", "10:41", "17_ROMERO",                107, "LP1912"),
    @(136, "08:47:19", "10:42", "17_ROMERO",                115, "LP1912"),
    @(137, "08:47:19", "10:43", "14_ABASTO",                116, "LP1912")
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:54:42"

# ---------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 08:54:42"
$ws3.Range("A3").Value = "Total filas: 25"

$sheet3Rows = @(
    @(27, "08:54:42", "09:11", "215D_LA PLATA",             17, "L6203"),
    @(28, "08:47:19", "09:12", "215D_LA PLATA",             25, "L6203"),
    @(29, "08:16:48", "10:02", "215B_LP-P MOR-40 Y 115",   106, "L6173"),
    @(30, "08:47:19", "10:03", "215B_LP-P MOR-40 Y 115",    76, "L6173")
)

foreach ($row in $sheet3Rows) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 2).Value = $row[2]
    $ws3.Cells.Item($r, 3).Value = $row[3]
    $ws3.Cells.Item($r, 4).Value = $row[4]
    $ws3.Cells.Item($r, 5).Value = $row[5]
}
